$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.754521
$ws.Range("H2").Value = 2.263563
$ws.Range("I2").Value = 0.2768403531129761
$ws.Range("J2").Value = 0.2768403531129761
$ws.Range("M2").Value = 41.60833666666667
$ws.Range("N2").Value = 124.82501
$ws.Range("O2").Value = 0.5886423873735626
$ws.Range("P2").Value = 0.5886423873735626
$ws.Range("Q2").Value = 31.39436379007
$ws.Range("R2").Value = 282.54927411063
$ws.Range("S2").Value = 0.1629599663777623
$ws.Range("T2").Value = 0.1629599663777623

# Row 3
$ws.Range("G3").Value = 0.754521
$ws.Range("H3").Value = 2.263563
$ws.Range("I3").Value = 0.2768403531129761
$ws.Range("J3").Value = 0.2768403531129761
$ws.Range("O3").Value = 0.2506218293658061
$ws.Range("P3").Value = 0.2506218293658061
$ws.Range("Q3").Value = 13.366541475121
$ws.Range("R3").Value = 120.298873276089
$ws.Range("S3").Value = 0.06938223573944981
$ws.Range("T3").Value = 0.06938223573944981

# Row 4
$ws.Range("G4").Value = 0.754521
$ws.Range("H4").Value = 2.263563
$ws.Range("I4").Value = 0.2768403531129761
$ws.Range("J4").Value = 0.2768403531129761
$ws.Range("M4").Value = 6.935318333333334
$ws.Range("N4").Value = 20.805955
$ws.Range("O4").Value = 0.09811549001908282
$ws.Range("P4").Value = 0.09811549001908279
$ws.Range("Q4").Value = 5.232843324185001
$ws.Range("R4").Value = 47.095589917665
$ws.Range("S4").Value = 0.02716232690273557
$ws.Range("T4").Value = 0.02716232690273556

# Row 5
$ws.Range("G5").Value = 0.754521
$ws.Range("H5").Value = 2.263563
$ws.Range("I5").Value = 0.2768403531129761
$ws.Range("J5").Value = 0.2768403531129761
$ws.Range("M5").Value = 4.426331333333334
$ws.Range("N5").Value = 13.278994
$ws.Range("O5").Value = 0.06262029324154843
$ws.Range("P5").Value = 0.06262029324154841
$ws.Range("Q5").Value = 3.339759943958
$ws.Range("R5").Value = 30.057839495622
$ws.Range("S5").Value = 0.01733582409302838
$ws.Range("T5").Value = 0.01733582409302837

# Row 6
$ws.Range("G6").Value = 0.9731926666666667
$ws.Range("I6").Value = 0.3570728998754956
$ws.Range("J6").Value = 0.3570728998754956
$ws.Range("M6").Value = 41.60833666666667
$ws.Range("N6").Value = 124.82501
$ws.Range("O6").Value = 0.5886423873735626
$ws.Range("P6").Value = 0.5886423873735626
$ws.Range("Q6").Value = 40.49292811619778
$ws.Range("R6").Value = 364.43635304578
$ws.Range("S6").Value = 0.2101882442491128
$ws.Range("T6").Value = 0.2101882442491128

# Row 7
$ws.Range("G7").Value = 0.9731926666666667
$ws.Range("I7").Value = 0.3570728998754956
$ws.Range("J7").Value = 0.3570728998754956
$ws.Range("O7").Value = 0.2506218293658061
$ws.Range("P7").Value = 0.2506218293658061
$ws.Range("S7").Value = 0.08949026338375003
$ws.Range("T7").Value = 0.08949026338375003

# Row 8
$ws.Range("G8").Value = 0.9731926666666667
$ws.Range("I8").Value = 0.3570728998754956
$ws.Range("J8").Value = 0.3570728998754956
$ws.Range("M8").Value = 6.935318333333334
$ws.Range("N8").Value = 20.805955
$ws.Range("O8").Value = 0.09811549001908282
$ws.Range("P8").Value = 0.09811549001908279
$ws.Range("Q8").Value = 6.74940094299889
$ws.Range("R8").Value = 60.74460848699
$ws.Range("S8").Value = 0.03503438254381915
$ws.Range("T8").Value = 0.03503438254381914

# Row 9
$ws.Range("G9").Value = 0.9731926666666667
$ws.Range("I9").Value = 0.3570728998754956
$ws.Range("J9").Value = 0.3570728998754956
$ws.Range("M9").Value = 4.426331333333334
$ws.Range("N9").Value = 13.278994
$ws.Range("O9").Value = 0.06262029324154843
$ws.Range("P9").Value = 0.06262029324154841
$ws.Range("Q9").Value = 4.307673193836889
$ws.Range("R9").Value = 38.769058744532
$ws.Range("S9").Value = 0.0223600096988136
$ws.Range("T9").Value = 0.02236000969881359

# Row 10
$ws.Range("G10").Value = 0.7824410000000001
$ws.Range("H10").Value = 2.347323
$ws.Range("I10").Value = 0.2870844452706686
$ws.Range("J10").Value = 0.2870844452706686
$ws.Range("M10").Value = 41.60833666666667
$ws.Range("N10").Value = 124.82501
$ws.Range("O10").Value = 0.5886423873735626
$ws.Range("P10").Value = 0.5886423873735626
$ws.Range("Q10").Value = 32.55606854980334
$ws.Range("R10").Value = 293.00461694823
$ws.Range("S10").Value = 0.1689900732419412
$ws.Range("T10").Value = 0.1689900732419412

# Row 11
$ws.Range("G11").Value = 0.7824410000000001
$ws.Range("H11").Value = 2.347323
$ws.Range("I11").Value = 0.2870844452706686
$ws.Range("J11").Value = 0.2870844452706686
$ws.Range("O11").Value = 0.2506218293658061
$ws.Range("P11").Value = 0.2506218293658061
$ws.Range("Q11").Value = 13.86115174837433
$ws.Range("R11").Value = 124.750365735369
$ws.Range("S11").Value = 0.07194962885620261
$ws.Range("T11").Value = 0.07194962885620262

# Row 12
$ws.Range("G12").Value = 0.7824410000000001
$ws.Range("H12").Value = 2.347323
$ws.Range("I12").Value = 0.2870844452706686
$ws.Range("J12").Value = 0.2870844452706686
$ws.Range("M12").Value = 6.935318333333334
$ws.Range("N12").Value = 20.805955
$ws.Range("O12").Value = 0.09811549001908282
$ws.Range("P12").Value = 0.09811549001908279
$ws.Range("Q12").Value = 5.426477412051668
$ws.Range("R12").Value = 48.83829670846501
$ws.Range("S12").Value = 0.02816743102458821
$ws.Range("T12").Value = 0.02816743102458821

# Row 13
$ws.Range("G13").Value = 0.7824410000000001
$ws.Range("H13").Value = 2.347323
$ws.Range("I13").Value = 0.2870844452706686
$ws.Range("J13").Value = 0.2870844452706686
$ws.Range("M13").Value = 4.426331333333334
$ws.Range("N13").Value = 13.278994
$ws.Range("O13").Value = 0.06262029324154843
$ws.Range("P13").Value = 0.06262029324154841
$ws.Range("Q13").Value = 3.463343114784667
$ws.Range("R13").Value = 31.170088033062
$ws.Range("S13").Value = 0.01797731214793653
$ws.Range("T13").Value = 0.01797731214793653

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.2153186666666667
$ws.Range("H14").Value = 0.645956
$ws.Range("I14").Value = 0.07900230174085969
$ws.Range("J14").Value = 0.07900230174085969
$ws.Range("M14").Value = 41.60833666666667
$ws.Range("N14").Value = 124.82501
$ws.Range("O14").Value = 0.5886423873735626
$ws.Range("P14").Value = 0.5886423873735626
$ws.Range("Q14").Value = 8.959051573284444
$ws.Range("R14").Value = 80.63146415956
$ws.Range("S14").Value = 0.04650410350474621
$ws.Range("T14").Value = 0.04650410350474621

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.2153186666666667
$ws.Range("H15").Value = 0.645956
$ws.Range("I15").Value = 0.07900230174085969
$ws.Range("J15").Value = 0.07900230174085969
$ws.Range("O15").Value = 0.2506218293658061
$ws.Range("P15").Value = 0.2506218293658061
$ws.Range("Q15").Value = 3.814427813629778
$ws.Range("R15").Value = 34.329850322668
$ws.Range("S15").Value = 0.01979970138640367
$ws.Range("T15").Value = 0.01979970138640367

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.2153186666666667
$ws.Range("H16").Value = 0.645956
$ws.Range("I16").Value = 0.07900230174085969
$ws.Range("J16").Value = 0.07900230174085969
$ws.Range("M16").Value = 6.935318333333334
$ws.Range("N16").Value = 20.805955
$ws.Range("O16").Value = 0.09811549001908282
$ws.Range("P16").Value = 0.09811549001908279
$ws.Range("Q16").Value = 1.493303496442222
$ws.Range("R16").Value = 13.43973146798
$ws.Range("S16").Value = 0.007751349547939888
$ws.Range("T16").Value = 0.007751349547939886

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.2153186666666667
$ws.Range("H17").Value = 0.645956
$ws.Range("I17").Value = 0.07900230174085969
$ws.Range("J17").Value = 0.07900230174085969
$ws.Range("M17").Value = 4.426331333333334
$ws.Range("N17").Value = 13.278994
$ws.Range("O17").Value = 0.06262029324154843
$ws.Range("P17").Value = 0.06262029324154841
$ws.Range("Q17").Value = 0.9530717609182222
$ws.Range("R17").Value = 8.577645848264
$ws.Range("S17").Value = 0.004947147301769926
$ws.Range("T17").Value = 0.004947147301769925
